# Update the "last updated" timestamp string shown in row 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 16:20"

# Country names at these rows changed because the underlying dataset was
# re-sorted by case count; update the country label for each affected row
$ws.Cells.Item(73, 1).Value = "Kenia"
$ws.Cells.Item(74, 1).Value = "Sudan"
$ws.Cells.Item(75, 1).Value = "El Salvador"
$ws.Cells.Item(111, 1).Value = "Cuba"
$ws.Cells.Item(112, 1).Value = "Malaui"
$ws.Cells.Item(113, 1).Value = "Mali"
$ws.Cells.Item(114, 1).Value = "Libano"

# Refreshed case statistics: Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes
# Row 4
$ws.Cells.Item(4, 2).Value = 3481680
$ws.Cells.Item(4, 3).Value = 2197
$ws.Cells.Item(4, 4).Value = 1550471
$ws.Cells.Item(4, 5).Value = 1792918
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 44
$ws.Cells.Item(4, 8).Value = 138291
# Row 6
$ws.Cells.Item(6, 2).Value = 912570
$ws.Cells.Item(6, 3).Value = 4925
$ws.Cells.Item(6, 4).Value = 574225
$ws.Cells.Item(6, 5).Value = 314547
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 71
$ws.Cells.Item(6, 8).Value = 23798
# Row 17
$ws.Cells.Item(17, 2).Value = 237803
$ws.Cells.Item(17, 3).Value = 2692
$ws.Cells.Item(17, 4).Value = 177560
$ws.Cells.Item(17, 5).Value = 57960
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 40
$ws.Cells.Item(17, 8).Value = 2283
# Row 19
$ws.Cells.Item(19, 2).Value = 200451
$ws.Cells.Item(19, 3).Value = 15
$ws.Cells.Item(19, 4).Value = 185100
$ws.Cells.Item(19, 5).Value = 6212
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 9139
# Row 24
$ws.Cells.Item(24, 2).Value = 104533
$ws.Cells.Item(24, 3).Value = 517
$ws.Cells.Item(24, 4).Value = 101160
$ws.Cells.Item(24, 5).Value = 3223
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 150
# Row 61
$ws.Cells.Item(61, 2).Value = 19439
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 13033
$ws.Cells.Item(61, 5).Value = 5754
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 3
$ws.Cells.Item(61, 8).Value = 652
# Row 63
$ws.Cells.Item(63, 2).Value = 18983
$ws.Cells.Item(63, 3).Value = 344
$ws.Cells.Item(63, 4).Value = 13991
$ws.Cells.Item(63, 5).Value = 4574
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 13
$ws.Cells.Item(63, 8).Value = 418
# Row 67
$ws.Cells.Item(67, 2).Value = 13993
$ws.Cells.Item(67, 3).Value = 402
$ws.Cells.Item(67, 4).Value = 8074
$ws.Cells.Item(67, 5).Value = 5854
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = 65
# Row 73
$ws.Cells.Item(73, 2).Value = 10791
$ws.Cells.Item(73, 3).Value = 497
$ws.Cells.Item(73, 4).Value = 3017
$ws.Cells.Item(73, 5).Value = 7572
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 5
$ws.Cells.Item(73, 8).Value = 202
# Row 74
$ws.Cells.Item(74, 2).Value = 10316
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 4).Value = 5403
$ws.Cells.Item(74, 5).Value = 4256
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 657
# Row 75
$ws.Cells.Item(75, 2).Value = 10303
$ws.Cells.Item(75, 3).Value = 325
$ws.Cells.Item(75, 4).Value = 5919
$ws.Cells.Item(75, 5).Value = 4106
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 11
$ws.Cells.Item(75, 8).Value = 278
# Row 78
$ws.Cells.Item(78, 2).Value = 8988
$ws.Cells.Item(78, 3).Value = 4
$ws.Cells.Item(78, 4).Value = 8138
$ws.Cells.Item(78, 5).Value = 597
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 253
# Row 96
$ws.Cells.Item(96, 2).Value = 4979
$ws.Cells.Item(96, 3).Value = 2
$ws.Cells.Item(96, 4).Value = 4743
$ws.Cells.Item(96, 5).Value = 180
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 56
# Row 110
$ws.Cells.Item(110, 2).Value = 2661
$ws.Cells.Item(110, 3).Value = 15
$ws.Cells.Item(110, 4).Value = 1988
$ws.Cells.Item(110, 5).Value = 662
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 11
# Row 111
$ws.Cells.Item(111, 2).Value = 2432
$ws.Cells.Item(111, 3).Value = 4
$ws.Cells.Item(111, 4).Value = 2275
$ws.Cells.Item(111, 5).Value = 70
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 87
# Row 112
$ws.Cells.Item(112, 2).Value = 2430
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 747
$ws.Cells.Item(112, 5).Value = 1644
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 39
# Row 113
$ws.Cells.Item(113, 2).Value = 2423
$ws.Cells.Item(113, 3).Value = 11
$ws.Cells.Item(113, 4).Value = 1748
$ws.Cells.Item(113, 5).Value = 554
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 121
# Row 114
$ws.Cells.Item(114, 2).Value = 2419
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 1423
$ws.Cells.Item(114, 5).Value = 960
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 36
# Row 160
$ws.Cells.Item(160, 2).Value = 439
$ws.Cells.Item(160, 3).Value = 22
$ws.Cells.Item(160, 4).Value = 138
$ws.Cells.Item(160, 5).Value = 280
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 2
$ws.Cells.Item(160, 8).Value = 21
